$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price / Volume(1h) columns are stored as literal text in the source feed
# (e.g. "294.13", "0.38%"), so force Text format per-cell before writing the new
# values; otherwise Excel would auto-coerce them into numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "294.13"
$ws.Range("E2").Value = "0.38%"
$ws.Range("D3").Value = "31.01"
$ws.Range("E3").Value = "0.01%"
$ws.Range("D4").Value = "4.939"
$ws.Range("E4").Value = "0.90%"
$ws.Range("D5").Value = "0.07338"
$ws.Range("E5").Value = "0.81%"
$ws.Range("D6").Value = "2.368"
$ws.Range("E6").Value = "32.70%"
$ws.Range("D7").Value = "7.718"
$ws.Range("E7").Value = "0.47%"
$ws.Range("E8").Value = "-0.64%"
$ws.Range("D9").Value = "0.9015"
$ws.Range("E9").Value = "-0.12%"
$ws.Range("D10").Value = "0.1680"
$ws.Range("E10").Value = "1.18%"
$ws.Range("D11").Value = "0.07959"
$ws.Range("D12").Value = "0.08028"
$ws.Range("E12").Value = "-0.57%"
$ws.Range("D13").Value = "0.03101"
$ws.Range("E13").Value = "1.30%"
$ws.Range("D14").Value = "0.1008"
$ws.Range("E14").Value = "0.67%"
$ws.Range("D15").Value = "0.001502"
$ws.Range("E15").Value = "-0.50%"
$ws.Range("D16").Value = "0.005783"
$ws.Range("E16").Value = "2.21%"
$ws.Range("D17").Value = "3.484"
$ws.Range("E17").Value = "0.71%"
$ws.Range("D18").Value = "2.074"
$ws.Range("E18").Value = "-1.56%"
$ws.Range("D19").Value = "0.3331"
$ws.Range("E19").Value = "1.03%"
$ws.Range("E20").Value = "-0.31%"
$ws.Range("D21").Value = "3.976"
$ws.Range("E21").Value = "-8.94%"
$ws.Range("D22").Value = "0.2096"
$ws.Range("E22").Value = "4.68%"
$ws.Range("D23").Value = "0.04543"
$ws.Range("E23").Value = "1.31%"
$ws.Range("D24").Value = "0.001210"
$ws.Range("E24").Value = "-0.54%"
$ws.Range("D25").Value = "0.004649"
$ws.Range("E25").Value = "15.16%"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").Value = "3.66%"
$ws.Range("D27").Value = "0.0003388"
$ws.Range("D39").Value = "0.01612"
$ws.Range("E39").Value = "-2.60%"
$ws.Range("D40").Value = "0.04457"
$ws.Range("E40").Value = "1.59%"
$ws.Range("D41").Value = "0.007366"
$ws.Range("E41").Value = "-0.96%"
$ws.Range("D42").Value = "0.1328"
$ws.Range("E42").Value = "0.67%"
$ws.Range("D43").Value = "0.008575"
$ws.Range("D44").Value = "0.001998"
$ws.Range("E44").Value = "-2.10%"
$ws.Range("D45").Value = "0.009493"
$ws.Range("E45").Value = "-6.68%"
$ws.Range("D46").Value = "0.00005914"
$ws.Range("E46").Value = "-2.26%"
$ws.Range("E47").Value = "-0.33%"
$ws.Range("E48").Value = "3.19%"
$ws.Range("D49").Value = "0.002893"
$ws.Range("E49").Value = "20.40%"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").Value = "-0.33%"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").Value = "-0.33%"
